# esta todo lo de notificaciones pero la grafica de geocercas se encima en mensaje

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# Add the new "Sensores" worksheet after Sheet1
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sensores"

# Header row
$ws2.Range("A1").Value = "nombre_usuario"
$ws2.Range("B1").Value = "ThinkPower  TL904D"
$ws2.Range("C1").Value = "AOVX GM100"

# Data rows
$ws2.Range("A2").Value = "Mario"
$ws2.Range("B2").Value = 1
$ws2.Range("C2").Value = 0

$ws2.Range("A3").Value = "luigi"
$ws2.Range("B3").Value = 1
$ws2.Range("C3").Value = 1

$ws2.Range("A4").Value = "peach"
$ws2.Range("B4").Value = 0
$ws2.Range("C4").Value = 1

# Copy header style from Sheet1's A1 (bold/border/center) to Sensores A1
[void]$ws1.Range("A1").Copy()
[void]$ws2.Range("A1").PasteSpecial(-4122)  # xlPasteFormats

# Selection / view tweaks on Sensores (becomes the active/tabSelected sheet)
[void]$ws2.Select()
[void]$ws2.Range("B1").Select()

# Selection tweak on Sheet1 (no longer tabSelected, selection is A1:A4 range w/o activeCell change)
[void]$ws1.Range("A1:A4").Select()

# Make Sensores the active sheet (matches activeTab="1" / tabSelected on sheet2)
[void]$ws2.Activate()
